# Weekly update: insert one new daily price record at row 123,
# pushing the existing rows 123:245 down to 124:246.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 123 (shifts 123:245 -> 124:246, inherits
# formatting - in particular column D's date style - from the row above).
$ws.Rows("123:123").Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A123").Value = 1
$ws.Range("B123").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C123").Value = "Arica y Parinacota"
$ws.Range("D123").Value = 44705
$ws.Range("E123").Value = 15
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100108
$ws.Range("H123").Value = "Tropicales y subtropicales"
$ws.Range("I123").Value = 100108006
$ws.Range("J123").Value = "Plátano"
$ws.Range("K123").Value = "Sin especificar"
$ws.Range("L123").Value = "Pintón"
$ws.Range("M123").Value = 120
$ws.Range("N123").Value = 15000
$ws.Range("O123").Value = 16000
$ws.Range("P123").Value = 15500
$ws.Range("Q123").Value = '$/caja 20 kilos'
$ws.Range("R123").Value = "Ecuador"
$ws.Range("S123").Value = 775
$ws.Range("T123").Value = 20
